# Insert a new data row before the current row 502, shifting all rows
# from 502..585 down to 503..586 (matches dimension change A1:R585 -> A1:R586).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(502).Insert()

# Populate the newly inserted row 502 with the new record.
$ws.Cells.Item(502, 1).Value = 8
$ws.Cells.Item(502, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(502, 3).Value = "Coquimbo"
$ws.Cells.Item(502, 4).Value = 45180
$ws.Cells.Item(502, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(502, 5).Value = 4
$ws.Cells.Item(502, 6).Value = 100112032
$ws.Cells.Item(502, 7).Value = "Zapallo italiano"
$ws.Cells.Item(502, 8).Value = "Sin especificar"
$ws.Cells.Item(502, 9).Value = "Primera"
$ws.Cells.Item(502, 10).Value = 400
$ws.Cells.Item(502, 11).Value = 14500
$ws.Cells.Item(502, 12).Value = 15000
$ws.Cells.Item(502, 13).Value = 14750
$ws.Cells.Item(502, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(502, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(502, 16).Value = 295
$ws.Cells.Item(502, 17).Value = 50
$ws.Cells.Item(502, 18).Value = "Hortaliza"
